# This script reassigns the Fecha/Calidad/Volumen/Precio min/Precio max/
# Precio promedio/Origen/Precio $/Kg values across the daily price rows
# (rows 2-41) of the sheet, per the weekly update described in the commit
# message "Fruta / hortaliza, semanal". Column A-C, E-K, Q, T are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 44354, "Primera", 150, 21000, 22000, 21500, "Perú", 1075)
    ,@(3, 44435, "Primera", 60, 25000, 25000, 25000, "Perú", 1250)
    ,@(4, 44620, "Primera", 60, 22000, 22000, 22000, "Perú", 1100)
    ,@(5, 44372, "Primera", 60, 20000, 21000, 20667, "Perú", 1033)
    ,@(6, 44333, "Primera", 30, 22000, 22000, 22000, "Perú", 1100)
    ,@(7, 44445, "Primera", 35, 20000, 20000, 20000, "Perú", 1000)
    ,@(8, 44166, "Primera", 120, 28000, 28000, 28000, "Perú", 1400)
    ,@(9, 44263, "Segunda", 150, 15000, 15000, 15000, "Perú", 750)
    ,@(10, 44417, "Primera", 30, 24000, 24000, 24000, "Perú", 1200)
    ,@(11, 44410, "Primera", 40, 25000, 25000, 25000, "Perú", 1250)
    ,@(12, 44357, "Primera", 200, 20000, 21000, 20500, "Perú", 1025)
    ,@(13, 44522, "Primera", 25, 30000, 30000, 30000, "Perú", 1500)
    ,@(14, 44363, "Primera", 150, 21000, 22000, 21500, "Perú", 1075)
    ,@(15, 44355, "Primera", 200, 20000, 21000, 20500, "Ecuador", 1025)
    ,@(16, 44830, "Primera", 200, 30000, 30000, 30000, "Perú", 1500)
    ,@(17, 44382, "Primera", 200, 19000, 20000, 19500, "Perú", 975)
    ,@(18, 44277, "Primera", 60, 24000, 24000, 24000, "Perú", 1200)
    ,@(19, 44326, "Primera", 40, 22000, 22000, 22000, "Perú", 1100)
    ,@(20, 44300, "Primera", 150, 19000, 20000, 19500, "Perú", 975)
    ,@(21, 44299, "Primera", 150, 19000, 20000, 19500, "Perú", 975)
    ,@(22, 44473, "Primera", 40, 24000, 24000, 24000, "Perú", 1200)
    ,@(23, 44298, "Primera", 240, 19000, 20000, 19500, "Perú", 975)
    ,@(24, 44452, "Primera", 35, 21000, 22000, 21429, "Perú", 1071)
    ,@(25, 44284, "Primera", 40, 23000, 23000, 23000, "Perú", 1150)
    ,@(26, 44613, "Primera", 60, 30000, 30000, 30000, "Perú", 1500)
    ,@(27, 44305, "Primera", 40, 24000, 24000, 24000, "Perú", 1200)
    ,@(28, 44312, "Primera", 50, 22000, 22000, 22000, "Perú", 1100)
    ,@(29, 44424, "Primera", 70, 24000, 25000, 24429, "Perú", 1221)
    ,@(30, 44396, "Primera", 45, 22000, 22000, 22000, "Perú", 1100)
    ,@(31, 44431, "Primera", 60, 25000, 25000, 25000, "Perú", 1250)
    ,@(32, 44760, "Primera", 300, 24000, 25000, 24500, "Perú", 1225)
    ,@(33, 44350, "Primera", 90, 21000, 22000, 21556, "Perú", 1078)
    ,@(34, 44438, "Primera", 25, 21000, 21000, 21000, "Perú", 1050)
    ,@(35, 44442, "Primera", 30, 22000, 22000, 22000, "Perú", 1100)
    ,@(36, 44302, "Primera", 100, 19000, 20000, 19500, "Perú", 975)
    ,@(37, 44365, "Primera", 150, 20000, 21000, 20500, "Perú", 1025)
    ,@(38, 44529, "Primera", 34, 28000, 28000, 28000, "Perú", 1400)
    ,@(39, 44356, "Primera", 100, 20000, 21000, 20500, "Perú", 1025)
    ,@(40, 44270, "Primera", 50, 24000, 24000, 24000, "Perú", 1200)
    ,@(41, 44165, "Primera", 300, 27000, 28000, 27500, "Perú", 1375)
)
foreach ($row in $data) {
    $r   = $row[0]
    $d   = $row[1]
    $l   = $row[2]
    $m   = $row[3]
    $n   = $row[4]
    $o   = $row[5]
    $p   = $row[6]
    $org = $row[7]
    $s   = $row[8]

    $ws.Range("D$r").Value = $d
    $ws.Range("L$r").Value = $l
    $ws.Range("M$r").Value = $m
    $ws.Range("N$r").Value = $n
    $ws.Range("O$r").Value = $o
    $ws.Range("P$r").Value = $p
    $ws.Range("R$r").Value = $org
    $ws.Range("S$r").Value = $s
}
